# Fruta / hortaliza, semanal
# Insert a new weekly record at the top of the "Feria Lagunitas de Puerto
# Montt - Pina" data block (row 446), pushing the existing rows 446:483
# down to 447:484.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 446:483 down one row, creating a fresh (blank) row 446.
$ws.Rows(446).Insert()

# Populate the new row 446 with this week's observation.
$ws.Range("A446").Value = 4
$ws.Range("B446").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C446").Value = "Los Lagos"
$ws.Range("D446").Value = 45223
$ws.Range("E446").Value = 10
$ws.Range("F446").Value = "Fruta"
$ws.Range("G446").Value = 100108
$ws.Range("H446").Value = "Tropicales y subtropicales"
$ws.Range("I446").Value = 100108005
$ws.Range("J446").Value = "Piña"
$ws.Range("K446").Value = "Caramelo"
$ws.Range("L446").Value = "Segunda"
$ws.Range("M446").Value = 100
$ws.Range("N446").Value = 25000
$ws.Range("O446").Value = 25000
$ws.Range("P446").Value = 25000
$ws.Range("Q446").Value = "$/caja 14 unidades"
$ws.Range("R446").Value = "Ecuador"
$ws.Range("S446").Value = 1786
$ws.Range("T446").Value = 14
